$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the existing 5-digit service ids with the hotel id (00105) so
# 10101 -> 1000510101, 10102 -> 1000510102, etc. (rows 2..11, column A)
$ws.Range("A2").Value = 1000510101
$ws.Range("A3").Value = 1000510102
$ws.Range("A4").Value = 1000510103
$ws.Range("A5").Value = 1000510104
$ws.Range("A6").Value = 1000510105
$ws.Range("A7").Value = 1000510106
$ws.Range("A8").Value = 1000510107
$ws.Range("A9").Value = 1000510108
$ws.Range("A10").Value = 1000510109
$ws.Range("A11").Value = 1000510110

# Column A needs to be a bit wider to fit the longer ids.
$ws.Columns("A").ColumnWidth = 10.35

# Selection moved from the single cell G4 to the whole id column (A2:A11).
[void]$ws.Range("A2:A11").Select()
